# Update the "datas" sheet: replace the old Field/Value/Url/hyperlink
# layout with the new Field/Value table (Region/District/Suburb_1/Suburb_2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old hyperlink (on B2) and its "Hyperlink" cell style before
# writing new content, so the rebuilt sheet carries no leftover link or
# link-styling.
$ws.Hyperlinks.Delete()
$ws.Range("B2").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

# Column A: field names (written first so they land first in the shared
# string table), column B: values.
$ws.Range("A1").Value = "Field"
$ws.Range("A2").Value = "Region"
$ws.Range("A3").Value = "District"
$ws.Range("A4").Value = "Suburb_1"
$ws.Range("A5").Value = "Suburb_2"

$ws.Range("B1").Value = "Value"
$ws.Range("B2").Value = "Wellington"
$ws.Range("B3").Value = "Wellington"
$ws.Range("B4").Value = "Brooklyn"
$ws.Range("B5").Value = "Kelburn"
